# Auto-generated Excel COM-interop script
# Applies numeric cell updates to the Sagittarius_Profits workbook sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the target diff.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 373.375
$ws.Range("I4").Value = 141.71428
$ws.Range("J4").Value = 1995
$ws.Range("K4").Value = 141.71428
$ws.Range("L4").Value = 1995
$ws.Range("M4").Value = -27.71428
$ws.Range("N4").Value = -2223
$ws.Range("H39").Value = 2816.6667
$ws.Range("I39").Value = 2816.6667
$ws.Range("K39").Value = 8450.000100000001
$ws.Range("M39").Value = -8154.000100000001
$ws.Range("H47").Value = 10500
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 10500
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 10500
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -12444
$ws.Range("H51").Value = 4462.125
$ws.Range("I51").Value = 3799.8333
$ws.Range("K51").Value = 3799.8333
$ws.Range("M51").Value = -3315.8333
$ws.Range("H80").Value = 3570.5334
$ws.Range("I80").Value = 3950.3333
$ws.Range("J80").Value = 3317.3333
$ws.Range("K80").Value = 11850.9999
$ws.Range("L80").Value = 9951.999899999999
$ws.Range("M80").Value = -10852.9999
$ws.Range("N80").Value = -11947.9999
$ws.Range("H83").Value = 3570.5334
$ws.Range("I83").Value = 3950.3333
$ws.Range("J83").Value = 3317.3333
$ws.Range("K83").Value = 35552.9997
$ws.Range("L83").Value = 29855.9997
$ws.Range("M83").Value = -30560.9997
$ws.Range("N83").Value = -39839.9997
$ws.Range("H106").Value = 4900.5
$ws.Range("I106").Value = 5357.7144
$ws.Range("J106").Value = 1700
$ws.Range("K106").Value = 5357.7144
$ws.Range("L106").Value = 1700
$ws.Range("M106").Value = -4726.7144
$ws.Range("N106").Value = -2962
$ws.Range("H116").Value = 6450
$ws.Range("I116").Value = 4000
$ws.Range("K116").Value = 4000
$ws.Range("M116").Value = -558

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1835436.9
$ws.Range("I32").Value = 2060937.6
$ws.Range("J32").Value = 877058.6
$ws.Range("K32").Value = 2060937.6
$ws.Range("L32").Value = 877058.6
$ws.Range("M32").Value = -2060650.6
$ws.Range("N32").Value = -877632.6
$ws.Range("H74").Value = 567.5
$ws.Range("I74").Value = 567.5
$ws.Range("K74").Value = 567.5
$ws.Range("M74").Value = 306.5
$ws.Range("H77").Value = 567.5
$ws.Range("I77").Value = 567.5
$ws.Range("K77").Value = 2837.5
$ws.Range("M77").Value = 1530.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 947.625
$ws.Range("J64").Value = 949.5
$ws.Range("L64").Value = 949.5
$ws.Range("N64").Value = -1399.5
$ws.Range("H67").Value = 947.625
$ws.Range("J67").Value = 949.5
$ws.Range("L67").Value = 949.5
$ws.Range("N67").Value = -2509.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 110.888885
$ws.Range("I7").Value = 74.75
$ws.Range("J7").Value = 139.8
$ws.Range("K7").Value = 74.75
$ws.Range("L7").Value = 139.8
$ws.Range("M7").Value = 38.25
$ws.Range("N7").Value = -365.8
$ws.Range("H58").Value = 1959.6
$ws.Range("I58").Value = 1849.5
$ws.Range("K58").Value = 1849.5
$ws.Range("M58").Value = -1646.5
$ws.Range("H62").Value = 2749.5
$ws.Range("I62").Value = 2499
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 2499
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -1875
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 2749.5
$ws.Range("I65").Value = 2499
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 12495
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -9375
$ws.Range("N65").Value = -21240
$ws.Range("H86").Value = 11134
$ws.Range("I86").Value = 14499.5
$ws.Range("K86").Value = 14499.5
$ws.Range("M86").Value = -13376.5
$ws.Range("H89").Value = 11134
$ws.Range("I89").Value = 14499.5
$ws.Range("K89").Value = 72497.5
$ws.Range("M89").Value = -66881.5
$ws.Range("H94").Value = 163545.72
$ws.Range("J94").Value = 5273.3335
$ws.Range("L94").Value = 5273.3335
$ws.Range("N94").Value = -6175.3335
$ws.Range("H107").Value = 998.73334
$ws.Range("I107").Value = 1159.75
$ws.Range("J107").Value = 814.7143
$ws.Range("K107").Value = 1159.75
$ws.Range("L107").Value = 814.7143
$ws.Range("M107").Value = 760.25
$ws.Range("N107").Value = -4654.7143
$ws.Range("H122").Value = 2483.5
$ws.Range("J122").Value = 2140.6667
$ws.Range("L122").Value = 6422.000100000001
$ws.Range("N122").Value = -11322.0001
$ws.Range("H134").Value = 2250
$ws.Range("I134").Value = 2000
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 6000
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -3465
$ws.Range("N134").Value = -12570
$ws.Range("H136").Value = 1959.6
$ws.Range("I136").Value = 1849.5
$ws.Range("K136").Value = 5548.5
$ws.Range("M136").Value = -2998.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1361.9333
$ws.Range("I51").Value = 485.7143
$ws.Range("K51").Value = 1457.1429
$ws.Range("M51").Value = -997.1428999999998
$ws.Range("H92").Value = 534.25
$ws.Range("J92").Value = 523.6667
$ws.Range("L92").Value = 1571.0001
$ws.Range("N92").Value = -4067.0001
$ws.Range("H93").Value = 20000
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 20000
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 60000
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -63744
$ws.Range("H98").Value = 1422.5
$ws.Range("J98").Value = 1582.3334
$ws.Range("L98").Value = 4747.0002
$ws.Range("N98").Value = -7743.0002
$ws.Range("H140").Value = 17432.8
$ws.Range("I140").Value = 4665.5
$ws.Range("J140").Value = 20624.625
$ws.Range("K140").Value = 13996.5
$ws.Range("L140").Value = 61873.875
$ws.Range("M140").Value = -8816.5
$ws.Range("N140").Value = -72233.875

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 2534.75
$ws.Range("I107").Value = 1313.6666
$ws.Range("J107").Value = 2941.7778
$ws.Range("K107").Value = 1313.6666
$ws.Range("L107").Value = 2941.7778
$ws.Range("M107").Value = 606.3334
$ws.Range("N107").Value = -6781.7778
$ws.Range("H113").Value = 875
$ws.Range("I113").Value = 838.8
$ws.Range("J113").Value = 911.2
$ws.Range("K113").Value = 838.8
$ws.Range("L113").Value = 911.2
$ws.Range("M113").Value = 1331.2
$ws.Range("N113").Value = -5251.2

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 833.1111
$ws.Range("I22").Value = 883
$ws.Range("J22").Value = 733.3333
$ws.Range("K22").Value = 883
$ws.Range("L22").Value = 733.3333
$ws.Range("M22").Value = -588
$ws.Range("N22").Value = -1323.3333
$ws.Range("H27").Value = 833.1111
$ws.Range("I27").Value = 883
$ws.Range("J27").Value = 733.3333
$ws.Range("K27").Value = 883
$ws.Range("L27").Value = 733.3333
$ws.Range("M27").Value = -776
$ws.Range("N27").Value = -947.3333
$ws.Range("H46").Value = 2937.111
$ws.Range("I46").Value = 1467
$ws.Range("K46").Value = 1467
$ws.Range("M46").Value = -1279
$ws.Range("H55").Value = 649.4706
$ws.Range("I55").Value = 399.33334
$ws.Range("J55").Value = 1249.8
$ws.Range("K55").Value = 399.33334
$ws.Range("L55").Value = 1249.8
$ws.Range("M55").Value = -226.33334
$ws.Range("N55").Value = -1595.8
$ws.Range("H61").Value = 3297.2
$ws.Range("I61").Value = 3297.2
$ws.Range("K61").Value = 3297.2
$ws.Range("M61").Value = -3095.2
$ws.Range("H82").Value = 3289.5
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 3289.5
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H113").Value = 3297.2
$ws.Range("I113").Value = 3297.2
$ws.Range("K113").Value = 3297.2
$ws.Range("M113").Value = -1127.2
$ws.Range("H136").Value = 1866.4166
$ws.Range("I136").Value = 831.0526
$ws.Range("K136").Value = 2493.1578
$ws.Range("M136").Value = 56.84220000000005

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 14132.429
$ws.Range("J62").Value = 9820.833000000001
$ws.Range("L62").Value = 9820.833000000001
$ws.Range("N62").Value = -11068.833
$ws.Range("H65").Value = 14132.429
$ws.Range("J65").Value = 9820.833000000001
$ws.Range("L65").Value = 49104.165
$ws.Range("N65").Value = -55344.165
$ws.Range("H81").Value = 476.22223
$ws.Range("I81").Value = 385.75
$ws.Range("J81").Value = 1200
$ws.Range("K81").Value = 771.5
$ws.Range("L81").Value = 2400
$ws.Range("M81").Value = 289.5
$ws.Range("N81").Value = -4522
$ws.Range("H84").Value = 476.22223
$ws.Range("I84").Value = 385.75
$ws.Range("J84").Value = 1200
$ws.Range("K84").Value = 3857.5
$ws.Range("L84").Value = 12000
$ws.Range("M84").Value = 1446.5
$ws.Range("N84").Value = -22608
$ws.Range("H107").Value = 518.8333
$ws.Range("J107").Value = 558
$ws.Range("L107").Value = 1674
$ws.Range("N107").Value = -5514
$ws.Range("H126").Value = 5166.6665
$ws.Range("I126").Value = 5166.6665
$ws.Range("K126").Value = 15499.9995
$ws.Range("M126").Value = -13029.9995
$ws.Range("H132").Value = 3181.2727
$ws.Range("I132").Value = 3442.7778
$ws.Range("J132").Value = 2004.5
$ws.Range("K132").Value = 10328.3334
$ws.Range("L132").Value = 6013.5
$ws.Range("M132").Value = -7798.3334
$ws.Range("N132").Value = -11073.5
